$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New named ranges (inserted alphabetically by Excel automatically)
# ---------------------------------------------------------------------------
$wb.Names.Add("B_round", "=גיליון1!`$A`$46:`$F`$51")
$wb.Names.Add("NSiB",    "=גיליון1!`$A`$53:`$F`$58")
$wb.Names.Add("TBS",     "=גיליון1!`$A`$60:`$F`$65")

# ---------------------------------------------------------------------------
# 2. Three new 6x6 data blocks below the existing "requests" table
# ---------------------------------------------------------------------------
$B_round = @(
    @(0,0,0,0,0,0),
    @(1,0,0,0,0,0),
    @(1,1,0,0,0,0),
    @(1,1,1,0,0,0),
    @(1,1,1,1,0,0),
    @(1,1,1,1,1,0)
)

$NSiB = @(
    @(0,1,2,3,4,5),
    @(5,0,1,2,3,4),
    @(4,5,0,1,2,3),
    @(3,4,5,0,1,2),
    @(2,3,4,5,0,1),
    @(1,2,3,4,5,0)
)

$TBS = @(
    @(0,2,4,6,8,10),
    @(10,0,2,4,6,8),
    @(8,10,0,2,4,6),
    @(6,8,10,0,2,4),
    @(4,6,8,10,0,2),
    @(2,4,6,8,10,0)
)

function Write-Block($ws, $startRow, $block) {
    for ($i = 0; $i -lt 6; $i++) {
        $row = $startRow + $i
        for ($j = 0; $j -lt 6; $j++) {
            $col = 1 + $j
            $ws.Cells.Item($row, $col).Value = $block[$i][$j]
        }
    }
}

Write-Block $ws 46 $B_round
Write-Block $ws 53 $NSiB
Write-Block $ws 60 $TBS

# ---------------------------------------------------------------------------
# 3. Thin box border drawn around each block
# ---------------------------------------------------------------------------
$ws.Range("A46:F51").BorderAround(1)
$ws.Range("A53:F58").BorderAround(1)
$ws.Range("A60:F65").BorderAround(1)

# ---------------------------------------------------------------------------
# 4. View state: scroll near the new tables and select C47
# ---------------------------------------------------------------------------
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 45
    $win.ScrollColumn = 1
} catch {
}

$ws.Range("C47").Select()
